$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update email addresses in column D
$ws.Range("D2").Value = "raj0011@mail.com"
$ws.Range("D3").Value = "taj0021@mail.com"
$ws.Range("D4").Value = "gani0031@bhai.com"

# Update the active cell selection shown in the sheet view
$ws.Range("P14").Select()
